# Updated symbol list on Thu Dec 15 06:58:02 UTC 2022 with GitHub Actions
# Applies the latest scraped coinranking.com price/volume snapshot to the
# "cryptos" sheet. Values are stored as plain text in the workbook, so each
# assignment uses a leading apostrophe to force Excel to keep them as text
# (matching the original inlineStr cells) instead of auto-converting to
# numbers and losing formatting such as trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

Set-TextValue "D2"  "264.83"
Set-TextValue "D5"  "0.06153"
Set-TextValue "D6"  "3.590"
Set-TextValue "D7"  "6.711"
Set-TextValue "D8"  "1.350"
Set-TextValue "D9"  "0.8296"
Set-TextValue "D10" "0.01360"
Set-TextValue "D11" "0.1579"
Set-TextValue "D12" "0.08251"
Set-TextValue "D14" "0.03167"
Set-TextValue "D15" "0.09249"
Set-TextValue "D16" "3.901"
Set-TextValue "D17" "0.001721"
Set-TextValue "D18" "0.04801"
Set-TextValue "D19" "0.006289"
Set-TextValue "D20" "0.005940"
Set-TextValue "D21" "0.001091"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.772"
Set-TextValue "D24" "2.300"
Set-TextValue "D40" "0.04639"

Set-TextValue "D43" "0.003132"
Set-TextValue "E43" "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.01038"
Set-TextValue "E44" "43LocalTradersLCT"

Set-TextValue "D45" "0.00006143"
Set-TextValue "D47" "0.7781"
Set-TextValue "D48" "0.2002"
